# The deck's Integral/"Red Violet" design (ppt/theme/theme1.xml, driven by the
# slide master) is being put back to the default Office color palette, which
# is the change a user makes from Design > Variants > Colors (or by re-picking
# the "Office" theme colors) without touching the font/format scheme, which
# already matches the Office defaults on both themes in this deck.
#
# Helper: turn a "RRGGBB" hex string into the BGR-packed Long that the
# PowerPoint object model's ColorFormat.RGB / ThemeColor.RGB expects.
function HexToRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

# ThemeColorScheme items are in fixed OOXML clrScheme order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink.
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $colorScheme.Item($i).RGB = HexToRgb $officeColors[$i - 1]
}
